# Atualizacao 16 nov 2020
# Adds one new monthly record (Mensal) and fifteen new daily records (Diario).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Mensal": append row 14
# ---------------------------------------------------------------
$wsMensal = $wb.Worksheets.Item("Mensal")

$mensalRow = 14
$wsMensal.Range("A$mensalRow").Value = 44150
$wsMensal.Range("B$mensalRow").Value = 95.28
$wsMensal.Range("C$mensalRow").Value = 110.6
$wsMensal.Range("D$mensalRow").Value = -13.85

# Copy the date number-format/style (bold, border, centered, yyyy-mm-dd) from
# the row above, same as every other row in column A.
$wsMensal.Range("A13").Copy()
$wsMensal.Range("A$mensalRow").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "Diario": append rows 368-382
# ---------------------------------------------------------------
$wsDiario = $wb.Worksheets.Item("Diario")

$diarioData = @(
    @(44136, 96.79000000000001, 110.6, -12.48),
    @(44137, 97.72, 110.6, -11.64),
    @(44138, 85.7, 110.6, -22.51),
    @(44139, 89.31999999999999, 110.6, -19.23),
    @(44140, 87.59999999999999, 110.6, -20.79),
    @(44141, 82.11, 110.6, -25.75),
    @(44142, 74.14, 110.6, -32.97),
    @(44143, 71.15000000000001, 110.6, -35.67),
    @(44144, 69.59999999999999, 110.6, -37.06),
    @(44145, 71.44, 110.6, -35.41),
    @(44146, 92.51000000000001, 110.6, -16.35),
    @(44147, 126.58, 110.6, 14.46),
    @(44148, 120, 110.6, 8.5),
    @(44149, 130.98, 110.6, 18.43),
    @(44150, 133.53, 110.6, 20.74)
)

$startRow = 368
for ($i = 0; $i -lt $diarioData.Count; $i++) {
    $r = $startRow + $i
    $vals = $diarioData[$i]

    $wsDiario.Range("A$r").Value = $vals[0]
    $wsDiario.Range("B$r").Value = $vals[1]
    $wsDiario.Range("C$r").Value = $vals[2]
    $wsDiario.Range("D$r").Value = $vals[3]

    $wsDiario.Range("A367").Copy()
    $wsDiario.Range("A$r").PasteSpecial(-4122)
}
